$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 106: Making Your Mark | Enchanted Palladium Ink
$ws.Range("H106").Value = 3338.1177
$ws.Range("I106").Value = 2799.8
$ws.Range("K106").Value = 2799.8
$ws.Range("M106").Value = -2168.8

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 5039.6
$ws.Range("I116").Value = 4132.6665
$ws.Range("K116").Value = 4132.6665
$ws.Range("M116").Value = -690.6665000000003

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4932.1387
$ws.Range("I138").Value = 2953.2
$ws.Range("J138").Value = 5984.766
$ws.Range("K138").Value = 8859.599999999999
$ws.Range("L138").Value = 17954.298
$ws.Range("M138").Value = -3719.599999999999
$ws.Range("N138").Value = -28234.298

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 24766.547
$ws.Range("I32").Value = 25540.822
$ws.Range("K32").Value = 25540.822
$ws.Range("M32").Value = -25253.822

# Row 43: They've Got Legs | Steel Sabatons
$ws.Range("H43").Value = 20609.5
$ws.Range("J43").Value = 22646
$ws.Range("L43").Value = 22646
$ws.Range("N43").Value = -23272

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1207.375
$ws.Range("I45").Value = 1207.375
$ws.Range("K45").Value = 1207.375
$ws.Range("M45").Value = -830.375

# Row 46: Get Me the Usual | Heavy Steel Flanchard
$ws.Range("H46").Value = 4081.5
$ws.Range("J46").Value = 3430.4
$ws.Range("L46").Value = 3430.4
$ws.Range("N46").Value = -4068.4

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 10995.952
$ws.Range("J61").Value = 10271
$ws.Range("L61").Value = 10271
$ws.Range("N61").Value = -10695

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 266026.44
$ws.Range("I74").Value = 287092.97
$ws.Range("J74").Value = 20250.334
$ws.Range("K74").Value = 287092.97
$ws.Range("L74").Value = 20250.334
$ws.Range("M74").Value = -286218.97
$ws.Range("N74").Value = -21998.334

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 266026.44
$ws.Range("I77").Value = 287092.97
$ws.Range("J77").Value = 20250.334
$ws.Range("K77").Value = 1435464.85
$ws.Range("L77").Value = 101251.67
$ws.Range("M77").Value = -1431096.85
$ws.Range("N77").Value = -109987.67

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 3495.6099
$ws.Range("I122").Value = 3221.2778
$ws.Range("K122").Value = 9663.8334
$ws.Range("M122").Value = -7213.8334

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 5166.109
$ws.Range("I132").Value = 3160.1538
$ws.Range("K132").Value = 9480.4614
$ws.Range("M132").Value = -6950.4614

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 10995.952
$ws.Range("J136").Value = 10271
$ws.Range("L136").Value = 30813
$ws.Range("N136").Value = -35913

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 7961.387
$ws.Range("I134").Value = 5273.6523
$ws.Range("J134").Value = 15688.625
$ws.Range("K134").Value = 15820.9569
$ws.Range("L134").Value = 47065.875
$ws.Range("M134").Value = -13285.9569
$ws.Range("N134").Value = -52135.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 23815700
$ws.Range("I31").Value = 100002504
$ws.Range("K31").Value = 100002504
$ws.Range("M31").Value = -100002209

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 23815700
$ws.Range("I34").Value = 100002504
$ws.Range("K34").Value = 100002504
$ws.Range("M34").Value = -100002302

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 1936.7059
$ws.Range("I122").Value = 1436
$ws.Range("K122").Value = 4308
$ws.Range("M122").Value = -1858

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 27516.637
$ws.Range("I132").Value = 1766.75
$ws.Range("K132").Value = 5300.25
$ws.Range("M132").Value = -2770.25

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2971.6667
$ws.Range("I134").Value = 2730.3726
$ws.Range("K134").Value = 8191.1178
$ws.Range("M134").Value = -5656.1178

$ws = $wb.Worksheets.Item("CUL")
# Row 56: Culture Club | Crowned Pie
$ws.Range("H56").Value = 7671.2856
$ws.Range("I56").Value = 7671.2856
$ws.Range("K56").Value = 7671.2856
$ws.Range("M56").Value = -7141.2856

# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 55560504
$ws.Range("I122").Value = 333357470
$ws.Range("J122").Value = 1111.0667
$ws.Range("K122").Value = 3000217230
$ws.Range("L122").Value = 9999.6003
$ws.Range("M122").Value = -3000214780
$ws.Range("N122").Value = -14899.6003

# Row 129: Comfort Food | Yakow Moussaka
$ws.Range("H129").Value = 29413808
$ws.Range("J129").Value = 71432720
$ws.Range("L129").Value = 214298160
$ws.Range("N129").Value = -214308160

$ws = $wb.Worksheets.Item("GSM")
# Row 43: Get the Green Stuff | Malachite Earrings
$ws.Range("H43").Value = 1214.2858
$ws.Range("I43").Value = 1214.2858
$ws.Range("K43").Value = 1214.2858
$ws.Range("M43").Value = -1063.2858

# Row 55: If You've Got It, Flaunt It | Peridot Earrings
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 58: The Big Red | Red Coral Necklace
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 1610.2
$ws.Range("I102").Value = 1188.5883
$ws.Range("K102").Value = 1188.5883
$ws.Range("M102").Value = 433.4117000000001

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 5520.4614
$ws.Range("I113").Value = 4275.1113
$ws.Range("J113").Value = 8322.5
$ws.Range("K113").Value = 4275.1113
$ws.Range("L113").Value = 8322.5
$ws.Range("M113").Value = -2105.1113
$ws.Range("N113").Value = -12662.5

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 1162.76
$ws.Range("I122").Value = 967.8570999999999
$ws.Range("K122").Value = 2903.5713
$ws.Range("M122").Value = -453.5712999999996

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 8007
$ws.Range("I132").Value = 2352.3333
$ws.Range("J132").Value = 11399.8
$ws.Range("K132").Value = 7056.999899999999
$ws.Range("L132").Value = 34199.39999999999
$ws.Range("M132").Value = -4526.999899999999
$ws.Range("N132").Value = -39259.39999999999

# Row 133: Pendulums of Our Own | Lar Pendulums
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 7966.278
$ws.Range("I40").Value = 8616.75
$ws.Range("K40").Value = 8616.75
$ws.Range("M40").Value = -8480.75

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 4483.1665
$ws.Range("I68").Value = 2979.8
$ws.Range("K68").Value = 2979.8
$ws.Range("M68").Value = -2230.8

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 4483.1665
$ws.Range("I71").Value = 2979.8
$ws.Range("K71").Value = 14899
$ws.Range("M71").Value = -11155

# Row 76: Dragoon Drop Rate | Dhalmelskin Breeches of Maiming
$ws.Range("H76").Value = 30262.334
$ws.Range("J76").Value = 30262.334
$ws.Range("L76").Value = 30262.334
$ws.Range("N76").Value = -30938.334

# Row 79: Exploiting the Adroit (L) | Dhalmelskin Breeches of Maiming
$ws.Range("H79").Value = 30262.334
$ws.Range("J79").Value = 30262.334
$ws.Range("L79").Value = 30262.334
$ws.Range("N79").Value = -32602.334

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 6609.355
$ws.Range("I132").Value = 5884.028
$ws.Range("J132").Value = 7613.654
$ws.Range("K132").Value = 17652.084
$ws.Range("L132").Value = 22840.962
$ws.Range("M132").Value = -15122.084
$ws.Range("N132").Value = -27900.962

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 6385.8936
$ws.Range("I136").Value = 4768.8823
$ws.Range("K136").Value = 14306.6469
$ws.Range("M136").Value = -11756.6469

$ws = $wb.Worksheets.Item("WVR")
# Row 49: A Leg Up on the Cold | Linen Tights
$ws.Range("H49").Value = 62983.5
$ws.Range("J49").Value = 62983.5
$ws.Range("L49").Value = 62983.5
$ws.Range("N49").Value = -63443.5

# Row 110: Suits You | Iridescent Acton of Aiming
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 3107.1296
$ws.Range("I132").Value = 2494.3809
$ws.Range("K132").Value = 7483.1427
$ws.Range("M132").Value = -4953.1427

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 3618.3794
$ws.Range("I136").Value = 2466.6
$ws.Range("K136").Value = 7399.799999999999
$ws.Range("M136").Value = -4849.799999999999
